# Regenerate save_data: update column G ("K") values for the leiter_mark sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of Excel row number -> new value for column G ("K")
$gValues = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 4
    6  = 2
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 3
    17 = 2
    18 = 3
    19 = 3
    20 = 1
    21 = 1
    22 = 3
    23 = 5
    24 = 1
    25 = 4
    26 = 5
    27 = 4
    28 = 2
    29 = 4
    30 = 3
    31 = 1
    32 = 2
    33 = 3
    34 = 0
    35 = 3
    36 = 4
    38 = 3
    39 = 3
    40 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
